$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 11:20"

# --- Row 20: Banglades (data refresh, no reorder) ---
$ws.Range("B20").Value = 181129
$ws.Range("C20").Value = 2686
$ws.Range("D20").Value = 88034
$ws.Range("E20").Value = 90790
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 2305

# --- Rows 29-30: Indonesia overtakes Irak ---
$ws.Range("A29").Value = "Indonesia"
$ws.Range("B29").Value = 74018
$ws.Range("C29").Value = 1671
$ws.Range("D29").Value = 34719
$ws.Range("E29").Value = 35764
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 66
$ws.Range("H29").Value = 3535

$ws.Range("A30").Value = "Irak"
$ws.Range("B30").Value = 72460
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 41380
$ws.Range("E30").Value = 28120
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 2960

# --- Rows 35-38: Oman & Filipinas overtake Emiratos Arabes Unidos & Kuwait ---
$ws.Range("A35").Value = "Oman"
$ws.Range("B35").Value = 54697
$ws.Range("C35").Value = 1083
$ws.Range("D35").Value = 35255
$ws.Range("E35").Value = 19194
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 248

$ws.Range("A36").Value = "Filipinas"
$ws.Range("B36").Value = 54222
$ws.Range("C36").Value = 1308
$ws.Range("D36").Value = 14037
$ws.Range("E36").Value = 38813
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 12
$ws.Range("H36").Value = 1372

$ws.Range("A37").Value = "Emiratos Arabes Unidos"
$ws.Range("B37").Value = 54050
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 43969
$ws.Range("E37").Value = 9751
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 330

$ws.Range("A38").Value = "Kuwait"
$ws.Range("B38").Value = 53580
$ws.Range("C38").Value = 0
$ws.Range("D38").Value = 43214
$ws.Range("E38").Value = 9983
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 383

# --- Row 46: Polonia (data refresh, no reorder) ---
$ws.Range("B46").Value = 37521
$ws.Range("C46").Value = 305
$ws.Range("D46").Value = 26635
$ws.Range("E46").Value = 9318
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 6
$ws.Range("H46").Value = 1568

# --- Row 61: Austria (data refresh, no reorder) ---
$ws.Range("B61").Value = 18783
$ws.Range("C61").Value = 74
$ws.Range("D61").Value = 16864
$ws.Range("E61").Value = 1213

# --- Row 83: Finlandia (data refresh, no reorder) ---
$ws.Range("B83").Value = 7291
$ws.Range("C83").Value = 12
$ws.Range("E83").Value = 162

# --- Row 121: Lituania (data refresh, no reorder) ---
$ws.Range("B121").Value = 1865
$ws.Range("C121").Value = 4
$ws.Range("D121").Value = 1579
$ws.Range("E121").Value = 207

# --- Row 123: Eslovenia (data refresh, no reorder) ---
$ws.Range("B123").Value = 1827
$ws.Range("C123").Value = 34
$ws.Range("E123").Value = 287

# --- Row 127: Hong Kong (data refresh, no reorder) ---
$ws.Range("B127").Value = 1433
$ws.Range("C127").Value = 29
$ws.Range("D127").Value = 1197
$ws.Range("E127").Value = 229

# --- Row 141: Uganda (data refresh, no reorder) ---
$ws.Range("B141").Value = 1013
$ws.Range("C141").Value = 7
$ws.Range("D141").Value = 952
$ws.Range("E141").Value = 61
